$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 21 (shifts the old rows 21-23 down to 22-24)
$ws.Rows("21").Insert() | Out-Null

# Populate the new row with the new skill enum entry
$ws.Range("G21").Value = "TIAN_SHAN_LIU_YANG_ZHANG"
$ws.Range("I21").Value = 6

# Match the resulting cursor position recorded in the saved workbook
$ws.Range("L19").Select() | Out-Null
